# Author typed two values on Sheet1 (A1 = "sa sa", A3 = "魚" with a
# furigana/phonetic guide reading that this engine doesn't round-trip),
# left the selection on A4, and the sheet was set up for printing on
# A4-portrait paper before the workbook was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value = "sa sa"
$ws.Range("A3").Value = "魚"

# Page setup (print dialog) — paper size 9 = A4, orientation 1 = portrait.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Leave the selection where the author left it before saving.
$ws.Range("A4").Select()
